$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 45 (existing rows 45-81 shift down to 46-82)
$ws.Rows.Item(45).Insert(-4121)

# Populate the newly inserted row 45 with the new weekly price record
$ws.Range("A45").Value = 10
$ws.Range("B45").Value = "Vega Modelo de Temuco"
$ws.Range("C45").Value = "La Araucanía"
$ws.Range("D45").Value = 45072
$ws.Range("E45").Value = 9
$ws.Range("F45").Value = 100112042
$ws.Range("G45").Value = "Locoto"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 50
$ws.Range("K45").Value = 4400
$ws.Range("L45").Value = 4400
$ws.Range("M45").Value = 4400
$ws.Range("N45").Value = "$/kilo"
$ws.Range("O45").Value = "Región de Arica y Parinacota"
$ws.Range("P45").Value = 4400
$ws.Range("Q45").Value = 1
$ws.Range("R45").Value = "Hortaliza"
